# fix: bonus penghasilan lainnya
#
# The template had three header columns: nip | kategori | nominal.
# The "kategori" column is dropped so the sheet becomes: nip | nominal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column B ("kategori"). A true column delete (not just a content
# clear) shifts column C ("nominal") left into B, and prunes the now-unused
# "kategori" entry out of the shared-strings table.
$ws.Columns.Item(2).Delete()

# The surviving first column (nip) gets a text ("@") number format, which
# is how the workbook ends up with a second cellXfs entry (numFmtId 49)
# applied to A1.
$ws.Range("A1").NumberFormat = "@"

# Move/restore the saved cursor position.
[void]$ws.Range("J17").Select()

# Touch page setup (portrait) so the worksheet carries an explicit
# <pageSetup> element, matching the saved print settings.
$ws.PageSetup.Orientation = 1
